# relatorio geral exportando em excel detalhado
#
# Turns the single-column "Requisições" report into a detailed, multi-column
# "Movimentações" report: renames the sheet, (re)writes the header row with
# 23 columns (A:W), widens/resizes the columns, and gives the last five
# headers (S:W) a gray header-fill style while keeping the bold+border style
# already used for the original header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet
$ws.Name = "Movimentações"

# 2. Header labels, written in the same order the original workbook's shared
#    string table uses them (so that sst indices line up with the source),
#    but each landing in its target column.
$ws.Range("A1").Value = "ID Movimentação Produto"
$ws.Range("B1").Value = "ID Orçamento"
$ws.Range("C1").Value = "ID Nota Fiscal"
$ws.Range("D1").Value = "ID Produto"
$ws.Range("G1").Value = "Forma de Pagamento"
$ws.Range("H1").Value = "Bandeira"
$ws.Range("P1").Value = "Nome do Produto"
$ws.Range("L1").Value = "Data de Cancelamento Nota Fiscal"
$ws.Range("K1").Value = "Data de Envio Nota Fiscal"
$ws.Range("J1").Value = "Status da Nota Fiscal"
$ws.Range("I1").Value = "Data de Criação Nota Fiscal"
$ws.Range("F1").Value = "Nome do Cliente"
$ws.Range("E1").Value = "Nome do Autor"
$ws.Range("M1").Value = "Chave de Acesso Nota Fiscal"
$ws.Range("N1").Value = "Protocolo de Envio Nota Fiscal"
$ws.Range("O1").Value = "Protocolo de Cancelamento Nota Fiscal"
$ws.Range("Q1").Value = "Valor da Unidade do Produto"
$ws.Range("R1").Value = "Estoque Atual do Produto"
$ws.Range("T1").Value = "Motivo da Movimentação"
$ws.Range("U1").Value = "Tipo de Movimentação"
$ws.Range("V1").Value = "Estoque Atual da Movimentação"
$ws.Range("W1").Value = "Valor Total da Movimentação"
$ws.Range("S1").Value = "Data da Movimentação"

# 3. Formatting: columns A:R reuse the bold+border header style that A1
#    already had; copy it across so no extra fonts/borders get created.
$ws.Range("A1").Copy()
$ws.Range("B1:R1").PasteSpecial(-4122)

# Columns S:W get the same bold+border style plus a gray fill, to set them
# apart as the "Movimentação" specific fields.
$ws.Range("S1:W1").PasteSpecial(-4122)
$ws.Range("S1:W1").Interior.Color = 12566463
$excel.CutCopyMode = 0

# 4. Column widths (header-fit-ish widths for the new layout)
$ws.Columns.Item(1).ColumnWidth = 23.736979166666668
$ws.Columns.Item(2).ColumnWidth = 12.451822916666666
$ws.Columns.Item(3).ColumnWidth = 12.022135416666666
$ws.Columns.Item(4).ColumnWidth = 9.592447916666666
$ws.Columns.Item(5).ColumnWidth = 29.877604166666668
$ws.Columns.Item(6).ColumnWidth = 29.877604166666668
$ws.Columns.Item(7).ColumnWidth = 19.166666666666668
$ws.Columns.Item(8).ColumnWidth = 10.592447916666666
$ws.Columns.Item(9).ColumnWidth = 24.166666666666668
$ws.Columns.Item(10).ColumnWidth = 18.451822916666668
$ws.Columns.Item(11).ColumnWidth = 22.592447916666668
$ws.Columns.Item(12).ColumnWidth = 30.736979166666668
$ws.Columns.Item(13).ColumnWidth = 25.451822916666668
$ws.Columns.Item(14).ColumnWidth = 27.307291666666668
$ws.Columns.Item(15).ColumnWidth = 35.451822916666664
$ws.Columns.Item(16).ColumnWidth = 29.877604166666668
$ws.Columns.Item(17).ColumnWidth = 26.307291666666668
$ws.Columns.Item(18).ColumnWidth = 23.166666666666668
$ws.Columns.Item(19).ColumnWidth = 23.166666666666668
$ws.Columns.Item(20).ColumnWidth = 23.307291666666668
$ws.Columns.Item(21).ColumnWidth = 20.877604166666668
$ws.Columns.Item(22).ColumnWidth = 29.451822916666668
$ws.Columns.Item(23).ColumnWidth = 26.592447916666668

# 5. View: scroll toward the right-hand columns and select S2, mirroring the
#    reviewer's viewport when the detailed report was last edited.
$ws.Range("S2").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 18
$win.ScrollRow = 1
